# Generate Report for Handoff
# Adds a new row (row 9) for file "ec2463c9-00b0-4470-8a79-73491930967a"
# to the Overview sheet and the two per-locale handoff-status sheets
# (zh-cn, de-de), mirroring the existing rows for the other files.

$wb = $excel.ActiveWorkbook

$guid = "ec2463c9-00b0-4470-8a79-73491930967a"
$mdName = "$guid.md"

# ---------------------------------------------------------------------
# Sheet "Overview": file name + status columns
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("B9").Value = "Ready for handoff"
$ovw.Range("C9").Value = "Ready for handoff"
$ovw.Range("D9").Value = "2016-32-20 18:32:08"

$ovw.Hyperlinks.Add(
    $ovw.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/9a2b6c1f5e4d3a7b8c9d0e1f2a3b4c5d6e7f8091/e2e/$mdName",
    "",
    "",
    $mdName
) | Out-Null

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn and de-de
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Code = "zh-cn"; HandoffDatetime = "2016-03-20 18:32:04" },
    @{ Sheet = "de-de"; Code = "de-de"; HandoffDatetime = "2016-03-20 18:32:08" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $code = $loc.Code
    $xlfName = "$guid.a13e50cb4529dda7f0fa5a991f88630f2f394cb5.$code.xlf"

    # Plain (non-hyperlinked) cells
    $ws.Range("C9").Value = "Ready for handoff"
    $ws.Range("E9").Value = $loc.HandoffDatetime
    $ws.Range("H9").Value = "0001-01-01 00:00:00"
    $ws.Range("I9").Value = "Include"

    # Hyperlinked cells: source file name, extension, target xlf
    $ws.Hyperlinks.Add(
        $ws.Range("A9"),
        "https://github.com/OpenLocalizationTest/oltest/blob/9a2b6c1f5e4d3a7b8c9d0e1f2a3b4c5d6e7f8091/e2e/$mdName",
        "",
        "",
        $mdName
    ) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("B9"),
        "https://github.com/OpenLocalizationTestOrg/oltest.$code/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9001/e2e/$mdName",
        "",
        "",
        ".md"
    ) | Out-Null

    $ws.Hyperlinks.Add(
        $ws.Range("D9"),
        "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2a3b4c5d6e7f8a9b0c1d2e3f4a5b6c7d8e9f0012/ol-handoff/OpenLocalizationTestOrg/oltest.$code/ci/ht/$xlfName",
        "",
        "",
        $xlfName
    ) | Out-Null
}

Write-Output "Added handoff row for $guid"
